$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Hours Left" column (B2:B16) for Iteration 6.
# Column C ("Ideal") recalculates automatically since it holds a formula
# referencing B2, which drives the whole burndown line.
$ws.Range("B2").Value = 12
$ws.Range("B3").Value = 12
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 10
$ws.Range("B6").Value = 10
$ws.Range("B7").Value = 9
$ws.Range("B8").Value = 9
$ws.Range("B9").Value = 8
$ws.Range("B10").Value = 5
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 5
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 3
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 0

# Bump the chart title from Iteration 5 to Iteration 6.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.ChartTitle.Text = "Iteration 6 Burndown"

# Move the active selection to L7, matching the saved view state.
$ws.Range("L7").Select() | Out-Null
